# Applies the Jan 26 2023 crypto price/volume refresh to Sheet1 columns D (Price) and E (Volume 1h).
# Values are written as text (apostrophe-prefixed) to match the source data which stores these
# as inline strings, not numbers/percentages; Style is reset to "Normal" afterwards so the
# forced-text formatting does not leave a visible quote-prefix style on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "307.69"),
    @("E2", "2.84%"),
    @("D3", "35.94"),
    @("E3", "2.04%"),
    @("D4", "5.067"),
    @("E4", "0.87%"),
    @("D5", "0.08123"),
    @("E5", "2.46%"),
    @("D6", "1.934"),
    @("E6", "2.22%"),
    @("D7", "4.161"),
    @("E7", "3.18%"),
    @("D8", "7.829"),
    @("E8", "0.58%"),
    @("D9", "0.9384"),
    @("E9", "1.52%"),
    @("D10", "0.1367"),
    @("E10", "-3.60%"),
    @("D11", "0.1915"),
    @("E11", "1.04%"),
    @("D12", "0.09223"),
    @("E12", "1.31%"),
    @("D13", "0.03507"),
    @("E13", "0.97%"),
    @("D14", "0.09896"),
    @("E14", "0.03%"),
    @("D15", "0.001455"),
    @("E15", "4.33%"),
    @("D16", "0.005828"),
    @("E16", "1.22%"),
    @("D17", "3.624"),
    @("E17", "3.32%"),
    @("D18", "2.941"),
    @("E18", "0.91%"),
    @("D19", "0.3451"),
    @("E19", "1.27%"),
    @("D20", "0.1345"),
    @("E20", "4.07%"),
    @("D21", "5.189"),
    @("E21", "2.66%"),
    @("E22", "5.11%"),
    @("D23", "0.04401"),
    @("E23", "-1.57%"),
    @("D24", "0.001234"),
    @("E24", "1.47%"),
    @("D25", "0.004770"),
    @("E25", "0.61%"),
    @("E26", "5.30%"),
    @("D27", "0.0003128"),
    @("E27", "3.91%"),
    @("D39", "0.02028"),
    @("E39", "7.59%"),
    @("D40", "0.05073"),
    @("E40", "8.16%"),
    @("D41", "0.01124"),
    @("E41", "15.77%"),
    @("D42", "0.007610"),
    @("E42", "3.23%"),
    @("D43", "0.1380"),
    @("E43", "4.77%"),
    @("E44", "2.44%"),
    @("D45", "0.01130"),
    @("E45", "20.95%"),
    @("D46", "0.00006313"),
    @("E46", "0.97%"),
    @("D47", "0.00000000750"),
    @("E47", "-0.32%"),
    @("D48", "65.22"),
    @("E48", "0.88%"),
    @("D49", "0.001190"),
    @("E49", "-28.49%"),
    @("D50", "0.00002100"),
    @("E50", "-0.32%"),
    @("D51", "0.0002000"),
    @("E51", "-0.32%")
)

foreach ($pair in $updates) {
    $cellRef = $pair[0]
    $newVal = $pair[1]
    $ws.Range($cellRef).Value = "'" + $newVal
    $ws.Range($cellRef).Style = "Normal"
}
